# edit.ps1
# Applies the "LinuxForHealth" rebrand edit to the StructureDefinition-care-gap-compliance-met
# workbook:
#   - Metadata sheet: URL, Version, Date, Publisher updated
#   - Elements sheet: the "Constraint(s)" cell for the root "Extension" row is cleared

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/care-gap-compliance-met"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 5 is the "Extension.url" element row; column Q is "Fixed Value", which
# mirrors the StructureDefinition's canonical URL shown on the Metadata sheet.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/care-gap-compliance-met"

# Row 2 is the "Extension" element row; column AI is "Constraint(s)".
# Clear the inherited ele-1/ext-1 constraint text that used to be shown here.
$elements.Range("AI2").Value = ""
